# Update the table style used by the three data tables in this deck
# (slides 14, 15 and 16) from the old custom "Table_0" style
# ({640A3823-F3E1-42A5-8B3B-85854ECCD087}) to PowerPoint's built-in
# "No Style, Table Grid" style ({5074AE63-8958-438D-B416-C99F4BDA33F2}).
#
# Table styles can't be reassigned by setting the Style property directly
# (PowerPoint raises "Table styles cannot be assigned through a property —
# call Table.ApplyStyle(...) instead"), so we use Table.ApplyStyle with the
# style's GUID, which is how the table style picker in the UI applies a
# style.

$targetStyleId = "{5074AE63-8958-438D-B416-C99F4BDA33F2}"

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($targetStyleId)
        }
    }
}
